$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly entered receiving data for rows 9 and 10
# Match the date number format already used by the cell above (J8)
$ws.Range("J9:J10").NumberFormat = "d-mmm"

$ws.Range("J9").Value = 44460
$ws.Range("K9").Value = 1746
$ws.Range("L9").Value = "36/20"
$ws.Range("M9").Value = 78
$ws.Range("N9").Value = 311

$ws.Range("J10").Value = 44461
$ws.Range("K10").Value = 1754
$ws.Range("L10").Value = "36/20"
$ws.Range("M10").Value = 78
$ws.Range("N10").Value = 153

# Recalculate so dependent formulas (O7, N31) update their cached values
$excel.Calculate()

# Update selection to match final author position
$ws.Range("N11").Select()
